$d = $word.ActiveDocument
$wns = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

# ---------------------------------------------------------------------------
# 1) Remove the old "FIRMA" paragraph that used to sit right before the page
#    break leading into the "CROQUIS DE UBICACION DEL DOMICILIO DEL
#    DEMANDADO" section.
# ---------------------------------------------------------------------------
$firmaParagraph = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text.TrimEnd("`r") -eq "FIRMA") {
        $firmaParagraph = $candidate
        break
    }
}
if ($firmaParagraph -ne $null) {
    $firmaParagraph.Range.Delete()
}

# ---------------------------------------------------------------------------
# 2) At the very end of the document (after the last picture, right before
#    the final section properties) append:
#       - three blank paragraphs
#       - a "List Bullet" paragraph: "Entrevista realizada por: David"
#       - a right aligned "FIRMA" paragraph (the one that moved down here)
# ---------------------------------------------------------------------------
function Append-ParagraphXml {
    param([string]$xml)
    $endPos = $d.Content.End
    $insertionRange = $d.Range($endPos - 1, $endPos - 1)
    $insertionRange.InsertParagraphAfter() | Out-Null
    $newParagraph = $d.Paragraphs.Item($d.Paragraphs.Count)
    $newParagraph.Range.InsertXML($xml) | Out-Null
}

Append-ParagraphXml "<w:p $wns/>"
Append-ParagraphXml "<w:p $wns/>"
Append-ParagraphXml "<w:p $wns/>"

Append-ParagraphXml ("<w:p $wns>" + `
    "<w:pPr><w:pStyle w:val='ListBullet'/></w:pPr>" + `
    "<w:r><w:rPr><w:b/><w:sz w:val='24'/></w:rPr>" + `
    "<w:t xml:space='preserve'>Entrevista realizada por: </w:t></w:r>" + `
    "<w:r><w:t>David</w:t></w:r></w:p>")

Append-ParagraphXml ("<w:p $wns>" + `
    "<w:pPr><w:jc w:val='right'/></w:pPr>" + `
    "<w:r><w:rPr><w:sz w:val='22'/></w:rPr><w:t>FIRMA</w:t></w:r></w:p>")
